$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9267526865005493
$ws.Range("B1").Value = 1.575770735740662
$ws.Range("C1").Value = 3.131964445114136
$ws.Range("D1").Value = 3.408442497253418
$ws.Range("E1").Value = 0.3659143447875977
